$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 653.58826
$ws.Range("I12").Value = 341.375
$ws.Range("K12").Value = 341.375
$ws.Range("M12").Value = -171.375

$ws.Range("H15").Value = 1371.8971
$ws.Range("I15").Value = 1371.8971
$ws.Range("K15").Value = 4115.6913
$ws.Range("M15").Value = -3946.6913

$ws.Range("H96").Value = 806590.9399999999
$ws.Range("I96").Value = 903.7778
$ws.Range("K96").Value = 2711.3334
$ws.Range("M96").Value = -1338.3334

$ws.Range("H111").Value = 3526
$ws.Range("I111").Value = 729.8
$ws.Range("J111").Value = 6322.2
$ws.Range("K111").Value = 2189.4
$ws.Range("L111").Value = 18966.6
$ws.Range("M111").Value = 877.6000000000004
$ws.Range("N111").Value = -25100.6

$ws.Range("H132").Value = 1801.8667
$ws.Range("I132").Value = 1839.9464
$ws.Range("J132").Value = 1268.75
$ws.Range("K132").Value = 5519.8392
$ws.Range("L132").Value = 3806.25
$ws.Range("M132").Value = -2989.8392
$ws.Range("N132").Value = -8866.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H18").Value = 1000000000
$ws.Range("J18").Value = 1000000000
$ws.Range("L18").Value = 1000000000
$ws.Range("N18").Value = -1000000644

$ws.Range("H32").Value = 9321.678
$ws.Range("I32").Value = 9168
$ws.Range("K32").Value = 9168
$ws.Range("M32").Value = -8881

$ws.Range("H46").Value = 27228.5
$ws.Range("J46").Value = 27228.5
$ws.Range("L46").Value = 27228.5
$ws.Range("N46").Value = -27866.5

$ws.Range("H53").Value = 250014940
$ws.Range("I53").Value = 9895
$ws.Range("K53").Value = 9895
$ws.Range("M53").Value = -9213

$ws.Range("H132").Value = 2706357.8
$ws.Range("I132").Value = 2905.1667
$ws.Range("K132").Value = 8715.500100000001
$ws.Range("M132").Value = -6185.500100000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H103").Value = 18219
$ws.Range("I103").Value = 0
$ws.Range("J103").Value = 18219
$ws.Range("K103").Value = 0
$ws.Range("L103").Value = 18219
$ws.Range("N103").Value = -20563
$ws.Range("M103").ClearContents()

$ws.Range("H105").Value = 332454.5
$ws.Range("I105").Value = 430237.06
$ws.Range("K105").Value = 430237.06
$ws.Range("M105").Value = -428490.06

$ws.Range("H107").Value = 2509.0967
$ws.Range("I107").Value = 2669.28
$ws.Range("J107").Value = 1841.6666
$ws.Range("K107").Value = 2669.28
$ws.Range("L107").Value = 1841.6666
$ws.Range("M107").Value = -749.2800000000002
$ws.Range("N107").Value = -5681.6666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 6599.4
$ws.Range("I58").Value = 2499.5
$ws.Range("K58").Value = 2499.5
$ws.Range("M58").Value = -2296.5

$ws.Range("H132").Value = 2836.1765
$ws.Range("I132").Value = 2863.1538
$ws.Range("K132").Value = 8589.4614
$ws.Range("M132").Value = -6059.4614

$ws.Range("H136").Value = 6599.4
$ws.Range("I136").Value = 2499.5
$ws.Range("K136").Value = 7498.5
$ws.Range("M136").Value = -4948.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 947.9643
$ws.Range("I5").Value = 699
$ws.Range("J5").Value = 1163.7333
$ws.Range("K5").Value = 2097
$ws.Range("L5").Value = 3491.199900000001
$ws.Range("M5").Value = -1985
$ws.Range("N5").Value = -3715.199900000001

$ws.Range("H135").Value = 947.9643
$ws.Range("I135").Value = 699
$ws.Range("J135").Value = 1163.7333
$ws.Range("K135").Value = 6291
$ws.Range("L135").Value = 10473.5997
$ws.Range("M135").Value = -3756
$ws.Range("N135").Value = -15543.5997

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H51").Value = 114612.5
$ws.Range("J51").Value = 114612.5
$ws.Range("L51").Value = 114612.5
$ws.Range("N51").Value = -115630.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("N18").ClearContents()

$ws.Range("H68").Value = 2456544.5
$ws.Range("I68").Value = 3792114.5
$ws.Range("J68").Value = 7999.6665
$ws.Range("K68").Value = 3792114.5
$ws.Range("L68").Value = 7999.6665
$ws.Range("M68").Value = -3791365.5
$ws.Range("N68").Value = -9497.666499999999

$ws.Range("H71").Value = 2456544.5
$ws.Range("I71").Value = 3792114.5
$ws.Range("J71").Value = 7999.6665
$ws.Range("K71").Value = 18960572.5
$ws.Range("L71").Value = 39998.3325
$ws.Range("M71").Value = -18956828.5
$ws.Range("N71").Value = -47486.3325

$ws.Range("H132").Value = 2959.9792
$ws.Range("I132").Value = 2180.2068
$ws.Range("J132").Value = 4150.1577
$ws.Range("K132").Value = 6540.6204
$ws.Range("L132").Value = 12450.4731
$ws.Range("M132").Value = -4010.6204
$ws.Range("N132").Value = -17510.4731

$ws.Range("H136").Value = 4377.185
$ws.Range("I136").Value = 3189.4285
$ws.Range("J136").Value = 8534.333000000001
$ws.Range("K136").Value = 9568.2855
$ws.Range("L136").Value = 25602.999
$ws.Range("M136").Value = -7018.2855
$ws.Range("N136").Value = -30702.999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H16").Value = 68887.60000000001
$ws.Range("J16").Value = 68887.60000000001
$ws.Range("L16").Value = 68887.60000000001
$ws.Range("N16").Value = -69471.60000000001

$ws.Range("H62").Value = 7835.2173
$ws.Range("I62").Value = 4915.5
$ws.Range("K62").Value = 4915.5
$ws.Range("M62").Value = -4291.5

$ws.Range("H65").Value = 7835.2173
$ws.Range("I65").Value = 4915.5
$ws.Range("K65").Value = 24577.5
$ws.Range("M65").Value = -21457.5

$ws.Range("H81").Value = 1467.2
$ws.Range("I81").Value = 1279.8462
$ws.Range("J81").Value = 2685
$ws.Range("K81").Value = 2559.6924
$ws.Range("L81").Value = 5370
$ws.Range("M81").Value = -1498.6924
$ws.Range("N81").Value = -7492

$ws.Range("H84").Value = 1467.2
$ws.Range("I84").Value = 1279.8462
$ws.Range("J84").Value = 2685
$ws.Range("K84").Value = 12798.462
$ws.Range("L84").Value = 26850
$ws.Range("M84").Value = -7494.462
$ws.Range("N84").Value = -37458

$ws.Range("H136").Value = 254187.03
$ws.Range("I136").Value = 4453.9697
$ws.Range("J136").Value = 1431500
$ws.Range("K136").Value = 13361.9091
$ws.Range("L136").Value = 4294500
$ws.Range("M136").Value = -10811.9091
$ws.Range("N136").Value = -4299600
